$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '68.040.28'
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").Value = '3.930.01'
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").Value = '476.43'
$ws.Range("E5").Value = '  +6.04%  '
$ws.Range("D6").Value = '146.85'
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '0.620'
$ws.Range("E7").Value = '  -1.03%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.728'
$ws.Range("E9").Value = '  -2.34%  '
$ws.Range("E10").Value = '  +4.39%  '
$ws.Range("D11").Value = '0.0000350'
$ws.Range("E11").Value = '  +7.02%  '
$ws.Range("D12").Value = '42.87'
$ws.Range("E12").Value = '  -2.39%  '
$ws.Range("D13").Value = '4.581.80'
$ws.Range("E13").Value = '  +3.14%  '
$ws.Range("D14").Value = '10.35'
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").Value = '4.015.71'
$ws.Range("E15").Value = '  +4.01%  '
$ws.Range("D16").Value = '14.70'
$ws.Range("E16").Value = '  -2.46%  '
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '19.87'
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("D19").Value = '1.13'
$ws.Range("E19").Value = '  -1.99%  '
$ws.Range("D20").Value = '68.320.40'
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("D21").Value = '437.66'
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("D22").Value = '14.40'
$ws.Range("E22").Value = '  -2.51%  '
$ws.Range("D23").Value = '3.31'
$ws.Range("E23").Value = '  +2.47%  '
$ws.Range("D24").Value = '87.40'
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("D25").Value = '3.63'
$ws.Range("E25").Value = '  +5.30%  '
$ws.Range("D26").Value = '38.15'
$ws.Range("E26").Value = '  +1.48%  '
$ws.Range("D27").Value = '10.26'
$ws.Range("E27").Value = '  +7.79%  '
$ws.Range("D28").Value = '10.18'
$ws.Range("E28").Value = '  +4.07%  '
$ws.Range("D29").Value = '727.80'
$ws.Range("E29").Value = '  -2.62%  '
$ws.Range("D30").Value = '13.28'
$ws.Range("E30").Value = '  -3.73%  '
$ws.Range("D31").Value = '0.128'
$ws.Range("E31").Value = '  -4.33%  '
$ws.Range("D32").Value = '2.80'
$ws.Range("E32").Value = '  +2.25%  '
$ws.Range("D33").Value = '41.98'
$ws.Range("E33").Value = '  -3.16%  '
$ws.Range("D34").Value = '0.0₃0872'
$ws.Range("E34").Value = '  +25.86%  '
$ws.Range("D35").Value = '59.51'
$ws.Range("E35").Value = '  +3.33%  '
$ws.Range("D36").Value = '0.151'
$ws.Range("E36").Value = '  -3.10%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '5.43'
$ws.Range("E37").Value = '  -1.66%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").Value = '0.0470'
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("D40").Value = '2.79'
$ws.Range("E40").Value = '  +13.09%  '
$ws.Range("D41").Value = '3.04'
$ws.Range("E41").Value = '  +4.45%  '
$ws.Range("D42").Value = '2.90'
$ws.Range("E42").Value = '  +9.09%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").Value = '0.344'
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '0.141'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '2.16'
$ws.Range("E46").Value = '  +0.34%  '
$ws.Range("B47").Value = 'LidoDAOToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D47").Value = '3.39'
$ws.Range("E47").Value = '  -2.68%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '146.50'
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '3.17'
$ws.Range("E49").Value = '  -1.56%  '
$ws.Range("D50").Value = '2.87'
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("D51").Value = '24.63'
$ws.Range("E51").Value = '  -2.64%  '
